$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")

# New header cell O1 ("Correction") — mirror the existing header style (N1)
# used for the other column headers (bold, centered, bordered).
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Cells.Item(1, 15).Value = "Correction"

# Data rows 2-12: column N ("Event") gets the literal text "nan", and a new
# (empty) column O is introduced alongside it.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
    # Materialize an empty cell at O<r> (matches the source which adds an
    # empty placeholder cell for every data row in the new column) without
    # pulling in any new number-format/style records.
    $ws.Cells.Item($r, 15).Font.Bold = $false
}
